$d = $word.ActiveDocument
$r = $d.Content
$r.InsertAfter("HELLO_WORLD_TEST")
